$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.082.80"
$ws.Range("E2").Value = "  +1.44%  "
$ws.Range("D3").Value = "3.766.68"
$ws.Range("E3").Value = "  -0.67%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "'622.12"
$ws.Range("E5").Value = "  +3.64%  "
$ws.Range("E6").Value = "  -0.62%  "
$ws.Range("D7").Value = "3.765.07"
$ws.Range("E7").Value = "  -0.69%  "
$ws.Range("E8").Value = "  -0.13%  "
$ws.Range("E9").Value = "  +0.50%  "
$ws.Range("E10").Value = "  +0.88%  "
$ws.Range("D11").Value = "'0.450"
$ws.Range("E11").Value = "  +0.19%  "
$ws.Range("D12").Value = "'6.61"
$ws.Range("E12").Value = "  +1.31%  "
$ws.Range("D13").Value = "'0.0000246"
$ws.Range("E13").Value = "  -0.60%  "
$ws.Range("D14").Value = "'35.31"
$ws.Range("E14").Value = "  -0.95%  "
$ws.Range("D15").Value = "4.403.88"
$ws.Range("E15").Value = "  -0.52%  "
$ws.Range("D16").Value = "3.814.45"
$ws.Range("E16").Value = "  +1.19%  "
$ws.Range("D17").Value = "69.105.92"
$ws.Range("E17").Value = "  +1.52%  "
$ws.Range("D18").Value = "'17.69"
$ws.Range("E18").Value = "  -3.17%  "
$ws.Range("D19").Value = "'7.07"
$ws.Range("E19").Value = "  +0.25%  "
$ws.Range("E20").Value = "  -1.28%  "
$ws.Range("D21").Value = "'467.49"
$ws.Range("E21").Value = "  +1.52%  "
$ws.Range("D22").Value = "'9.58"
$ws.Range("E22").Value = "  -0.88%  "
$ws.Range("D23").Value = "'0.700"
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("E24").Value = "  +1.92%  "
$ws.Range("D25").Value = "'82.87"
$ws.Range("E25").Value = "  -0.10%  "
$ws.Range("D26").Value = "'12.00"
$ws.Range("E26").Value = "  +0.27%  "
$ws.Range("E27").Value = "  +2.19%  "
$ws.Range("D28").Value = "'0.999"
$ws.Range("E28").Value = "  -0.09%  "
$ws.Range("D29").Value = "'9.96"
$ws.Range("E29").Value = "  -0.12%  "
$ws.Range("D30").Value = "3.914.90"
$ws.Range("E30").Value = "  -0.64%  "
$ws.Range("D31").Value = "'2.65"
$ws.Range("E31").Value = "  +0.73%  "
$ws.Range("D32").Value = "'2.24"
$ws.Range("E32").Value = "  +0.14%  "
$ws.Range("D33").Value = "'7.26"
$ws.Range("E33").Value = "  -0.44%  "
$ws.Range("D34").Value = "'28.82"
$ws.Range("E34").Value = "  -1.44%  "
$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = "  -0.06%  "
$ws.Range("D36").Value = "3.718.45"
$ws.Range("E36").Value = "  -0.56%  "
$ws.Range("D37").Value = "'8.94"
$ws.Range("E37").Value = "  -0.57%  "
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").Value = "'0.154"
$ws.Range("E38").Value = "  +10.30%  "
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").Value = "'0.102"
$ws.Range("E39").Value = "  +2.19%  "
$ws.Range("D40").Value = "'3.32"
$ws.Range("E40").Value = "  +0.39%  "
$ws.Range("D41").Value = "'5.76"
$ws.Range("E41").Value = "  -0.86%  "
$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").Value = "'1.00"
$ws.Range("E42").Value = "  +0.24%  "
$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D43").Value = "'0.964"
$ws.Range("E43").Value = "  -2.24%  "
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("D45").Value = "'0.299"
$ws.Range("E45").Value = "  +0.16%  "
$ws.Range("D46").Value = "'153.67"
$ws.Range("E46").Value = "  +1.15%  "
$ws.Range("D47").Value = "'43.01"
$ws.Range("E47").Value = "  -0.88%  "
$ws.Range("B48").Value = "OKB"
$ws.Range("C48").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D48").Value = "'46.68"
$ws.Range("E48").Value = "  -1.65%  "
$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").Value = "'1.91"
$ws.Range("E49").Value = "  +2.38%  "
$ws.Range("D50").Value = "'8.38"
$ws.Range("E50").Value = "  +0.42%  "
$ws.Range("E51").Value = "  +0.10%  "
